$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $row = $ws.Rows.Item(2)
    $origHeight = $row.RowHeight

    $val = $cell.Value()

    $newVal = $val.Replace("Date:                Sun, 05 Jan 2020", "Date:                Wed, 08 Jan 2020")
    $newVal = $newVal.Replace("Time:                        21:22:36", "Time:                        19:07:41")
    $newVal = $newVal.Replace("Time:                        21:22:37", "Time:                        19:07:41")

    $cell.Value = $newVal

    # Restore the row height, since assigning a new value to this
    # word-wrapped cell causes the engine to auto re-fit the row.
    $row.RowHeight = $origHeight
}
